$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.264.23"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "2.964.30"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("D5").Value = "'382.80"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").Value = "'103.32"
$ws.Range("E6").Value = "  -2.07%  "
$ws.Range("D7").Value = "'0.541"
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.590"
$ws.Range("D10").Value = "'36.69"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "'0.0842"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").Value = "3.427.10"
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("D14").Value = "'18.05"
$ws.Range("E14").Value = "  -2.97%  "
$ws.Range("D15").Value = "'7.47"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").Value = "2.956.50"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").Value = "'0.993"
$ws.Range("E17").Value = "  +3.82%  "
$ws.Range("D18").Value = "51.176.32"
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("E19").Value = "  -6.33%  "
$ws.Range("D20").Value = "'7.15"
$ws.Range("E20").Value = "  -3.92%  "
$ws.Range("D21").Value = "'12.63"
$ws.Range("E21").Value = "  -4.42%  "
$ws.Range("D22").Value = "0.0₃0957"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "'68.57"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'262.55"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").Value = "'2.93"
$ws.Range("E25").Value = "  +4.07%  "
$ws.Range("D26").Value = "'8.39"
$ws.Range("E26").Value = "  +12.73%  "
$ws.Range("D27").Value = "'7.84"
$ws.Range("E27").Value = "  +4.65%  "
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").Value = "'0.113"
$ws.Range("E29").Value = "  +8.53%  "
$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").Value = "'25.76"
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("D32").Value = "'9.84"
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("E33").Value = "  +5.22%  "
$ws.Range("D34").Value = "'34.01"
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").Value = "'2.06"
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "'50.46"
$ws.Range("E36").Value = "  -3.70%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").Value = "'3.00"
$ws.Range("E38").Value = "  -2.12%  "
$ws.Range("D39").Value = "'16.87"
$ws.Range("E39").Value = "  -2.89%  "
$ws.Range("D40").Value = "'2.57"
$ws.Range("E40").Value = "  -3.83%  "
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("E42").Value = "  -3.10%  "
$ws.Range("D43").Value = "'121.39"
$ws.Range("E43").Value = "  -2.69%  "
$ws.Range("D44").Value = "'21.47"
$ws.Range("E44").Value = "  -1.89%  "
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("E47").Value = "  +2.65%  "
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("D49").Value = "2.016.37"
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("D50").Value = "'0.0349"
$ws.Range("E50").Value = "  +6.59%  "
$ws.Range("D51").Value = "'2.13"
$ws.Range("E51").Value = "  +13.58%  "
